$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 3
$ws.Range("H3").Value = 2
$ws.Range("H2").Select()
